# Auto-generated Excel COM-interop script
# Applies scheduled-runner profit/price updates across multiple Leve sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC row 113
$ws.Range("H113").Value = 2493.45
$ws.Range("I113").Value = 2499.3333
$ws.Range("K113").Value = 2499.3333
$ws.Range("M113").Value = 754.6667000000002

# ALC row 137
$ws.Range("H137").Value = 3455602
$ws.Range("I137").Value = 5021.0454
$ws.Range("K137").Value = 15063.1362
$ws.Range("M137").Value = -12513.1362

# ALC row 138
$ws.Range("H138").Value = 4996.098
$ws.Range("I138").Value = 12499.777
$ws.Range("J138").Value = 3388.1667
$ws.Range("K138").Value = 37499.331
$ws.Range("L138").Value = 10164.5001
$ws.Range("M138").Value = -32359.331
$ws.Range("N138").Value = -20444.5001

$ws = $wb.Worksheets.Item("ARM")
# ARM row 32
$ws.Range("H32").Value = 2676.9434
$ws.Range("J32").Value = 7107.357
$ws.Range("L32").Value = 7107.357
$ws.Range("N32").Value = -7681.357

# ARM row 61
$ws.Range("H61").Value = 3210905.5
$ws.Range("I61").Value = 82526.86
$ws.Range("J61").Value = 11970365
$ws.Range("K61").Value = 82526.86
$ws.Range("L61").Value = 11970365
$ws.Range("M61").Value = -82314.86
$ws.Range("N61").Value = -11970789

# ARM row 74
$ws.Range("H74").Value = 620434.75
$ws.Range("I74").Value = 1063.1111
$ws.Range("K74").Value = 1063.1111
$ws.Range("M74").Value = -189.1111000000001

# ARM row 77
$ws.Range("H77").Value = 620434.75
$ws.Range("I77").Value = 1063.1111
$ws.Range("K77").Value = 5315.5555
$ws.Range("M77").Value = -947.5555000000004

# ARM row 88
$ws.Range("H88").Value = 1788.4
$ws.Range("I88").Value = 1950
$ws.Range("K88").Value = 1950
$ws.Range("M88").Value = -1544

# ARM row 91
$ws.Range("H91").Value = 1788.4
$ws.Range("I91").Value = 1950
$ws.Range("K91").Value = 1950
$ws.Range("M91").Value = -546

# ARM row 97
$ws.Range("H97").Value = 2432.8
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 2432.8
$ws.Range("K97").Value = 0
$ws.Range("L97").Value = 2432.8
$ws.Range("N97").Value = -3424.8
$ws.Range("M97").ClearContents()

# ARM row 122
$ws.Range("H122").Value = 1859.2
$ws.Range("I122").Value = 1798.8572
$ws.Range("K122").Value = 5396.571599999999
$ws.Range("M122").Value = -2946.571599999999

# ARM row 132
$ws.Range("H132").Value = 2860.24
$ws.Range("I132").Value = 2839.5334
$ws.Range("J132").Value = 2891.3
$ws.Range("K132").Value = 8518.600199999999
$ws.Range("L132").Value = 8673.900000000001
$ws.Range("M132").Value = -5988.600199999999
$ws.Range("N132").Value = -13733.9

# ARM row 136
$ws.Range("H136").Value = 3210905.5
$ws.Range("I136").Value = 82526.86
$ws.Range("J136").Value = 11970365
$ws.Range("K136").Value = 247580.58
$ws.Range("L136").Value = 35911095
$ws.Range("M136").Value = -245030.58
$ws.Range("N136").Value = -35916195

$ws = $wb.Worksheets.Item("BSM")
# BSM row 105
$ws.Range("H105").Value = 59499
$ws.Range("I105").Value = 99999
$ws.Range("J105").Value = 18999
$ws.Range("K105").Value = 99999
$ws.Range("L105").Value = 18999
$ws.Range("M105").Value = -98252
$ws.Range("N105").Value = -22493

$ws = $wb.Worksheets.Item("CRP")
# CRP row 22
$ws.Range("H22").Value = 1082
$ws.Range("I22").Value = 1299.75
$ws.Range("K22").Value = 1299.75
$ws.Range("M22").Value = -949.75

# CRP row 31
$ws.Range("H31").Value = 2871.1428
$ws.Range("I31").Value = 3866.2856
$ws.Range("J31").Value = 2373.5715
$ws.Range("K31").Value = 3866.2856
$ws.Range("L31").Value = 2373.5715
$ws.Range("M31").Value = -3571.2856
$ws.Range("N31").Value = -2963.5715

# CRP row 34
$ws.Range("H34").Value = 2871.1428
$ws.Range("I34").Value = 3866.2856
$ws.Range("J34").Value = 2373.5715
$ws.Range("K34").Value = 3866.2856
$ws.Range("L34").Value = 2373.5715
$ws.Range("M34").Value = -3664.2856
$ws.Range("N34").Value = -2777.5715

# CRP row 105
$ws.Range("H105").Value = 2488.182
$ws.Range("I105").Value = 1558.75
$ws.Range("K105").Value = 1558.75
$ws.Range("M105").Value = 188.25

# CRP row 122
$ws.Range("H122").Value = 3918.8262
$ws.Range("I122").Value = 3488.111
$ws.Range("K122").Value = 10464.333
$ws.Range("M122").Value = -8014.332999999999

$ws = $wb.Worksheets.Item("CUL")
# CUL row 51
$ws.Range("H51").Value = 3000
$ws.Range("J51").Value = 3000
$ws.Range("L51").Value = 9000
$ws.Range("N51").Value = -9920

# CUL row 104
$ws.Range("H104").Value = 7690.3335
$ws.Range("I104").Value = 2607.25
$ws.Range("J104").Value = 13499.571
$ws.Range("K104").Value = 7821.75
$ws.Range("L104").Value = 40498.713
$ws.Range("M104").Value = -5200.75
$ws.Range("N104").Value = -45740.713

# CUL row 105
$ws.Range("H105").Value = 16623.166
$ws.Range("J105").Value = 17947.8
$ws.Range("L105").Value = 53843.39999999999
$ws.Range("N105").Value = -59085.39999999999

# CUL row 106
$ws.Range("H106").Value = 13747.5
$ws.Range("I106").Value = 5000
$ws.Range("K106").Value = 15000
$ws.Range("M106").Value = -14054

# CUL row 122
$ws.Range("H122").Value = 6960778
$ws.Range("I122").Value = 33333604
$ws.Range("K122").Value = 300002436
$ws.Range("M122").Value = -299999986

$ws = $wb.Worksheets.Item("GSM")
# GSM row 2
$ws.Range("H2").Value = 5668.6665
$ws.Range("J2").Value = 10135.4
$ws.Range("L2").Value = 10135.4
$ws.Range("N2").Value = -10361.4

# GSM row 80
$ws.Range("H80").Value = 55576940
$ws.Range("I80").Value = 16114.75
$ws.Range("K80").Value = 16114.75
$ws.Range("M80").Value = -15116.75

# GSM row 83
$ws.Range("H83").Value = 55576940
$ws.Range("I83").Value = 16114.75
$ws.Range("K83").Value = 80573.75
$ws.Range("M83").Value = -75581.75

# GSM row 86
$ws.Range("H86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").ClearContents()

# GSM row 89
$ws.Range("H89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").ClearContents()

# GSM row 102
$ws.Range("H102").Value = 27778792
$ws.Range("I102").Value = 29412792
$ws.Range("K102").Value = 29412792
$ws.Range("M102").Value = -29411170

# GSM row 122
$ws.Range("H122").Value = 2214.2
$ws.Range("I122").Value = 2131.818
$ws.Range("J122").Value = 2440.75
$ws.Range("K122").Value = 6395.454000000001
$ws.Range("L122").Value = 7322.25
$ws.Range("M122").Value = -3945.454000000001
$ws.Range("N122").Value = -12222.25

# GSM row 132
$ws.Range("H132").Value = 10211135
$ws.Range("I132").Value = 4252.25
$ws.Range("J132").Value = 11986245
$ws.Range("K132").Value = 12756.75
$ws.Range("L132").Value = 35958735
$ws.Range("M132").Value = -10226.75
$ws.Range("N132").Value = -35963795

$ws = $wb.Worksheets.Item("LTW")
# LTW row 100
$ws.Range("H100").Value = 3157.5
$ws.Range("I100").Value = 2631.9092
$ws.Range("J100").Value = 3799.889
$ws.Range("K100").Value = 2631.9092
$ws.Range("L100").Value = 3799.889
$ws.Range("M100").Value = -2090.9092
$ws.Range("N100").Value = -4881.889

# LTW row 122
$ws.Range("H122").Value = 4145.125
$ws.Range("I122").Value = 3860.1667
$ws.Range("J122").Value = 5000
$ws.Range("K122").Value = 11580.5001
$ws.Range("L122").Value = 15000
$ws.Range("M122").Value = -9130.500100000001
$ws.Range("N122").Value = -19900

# LTW row 132
$ws.Range("H132").Value = 3495.5833
$ws.Range("I132").Value = 3581
$ws.Range("J132").Value = 3376
$ws.Range("K132").Value = 10743
$ws.Range("L132").Value = 10128
$ws.Range("M132").Value = -8213
$ws.Range("N132").Value = -15188

$ws = $wb.Worksheets.Item("WVR")
# WVR row 132
$ws.Range("H132").Value = 2917.889
$ws.Range("J132").Value = 3085.818
$ws.Range("L132").Value = 9257.454000000002
$ws.Range("N132").Value = -14317.454

Write-Output "Applied scheduled Sheets update across 8 worksheets."
